# Rebuild the document body to match the target edit:
# - Adds sz/szCs 22 (11pt) formatting to several heading/spacer paragraphs
# - Restructures the "2. ___ Lucas Lamounier Goncalves Duarte ___" run split
# - Replaces the placeholder bracketed instructional paragraphs with the
#   group's actual filled-in content (Aplicacao, objetivos, base de dados,
#   usuario, SGBD, linguagem de programacao)
# - Removes several now-unused blank paragraphs and the trailing
#   "[Assinatura ...]" placeholder paragraph
#
# Implemented via Range.InsertXML with a WordProcessingML package fragment
# covering the whole document body (sectPr is left untouched since it is
# outside $d.Content).

$d = $word.ActiveDocument

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:color w:val="002060"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:color w:val="002060"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:color w:val="002060"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">ANEXO I </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:color w:val="0070C0"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:color w:val="0070C0"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">Documento de Ciência do Grupo </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">Integrantes do Grupo: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>1. _________</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Kevin Vieira Pereira</w:t>
      </w:r>
      <w:r>
        <w:t>_______________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>2. ________</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>_</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve">Lucas </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Lamounier</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve"> Gonçalves Duarte</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>_________</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>3. _________</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Rodrigo Aparecido Silva Maia</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>______________</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">Descrição da Aplicação: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>plicação</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>a aplicação contará com uma interface simples, mas objetiva em que as opções de acessar relatórios, ou da criação do CRUD, sejam acessadas apenas pelos usuários designados.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Principais objetivos e público alvo</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> a aplicação tem como objetivo fornecer relatórios que permitam controle sobre o pagamento das diárias, por cada órgão</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. Tendo como público alvo o governo, apresentando os dados de uma maneira </w:t>
      </w:r>
      <w:r>
        <w:t>a qual</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> possa </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ajudá-lo </w:t>
      </w:r>
      <w:r>
        <w:t>a controlar de maneira mais eficiente o dinheiro destinado ao pagamento de diárias.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>B</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>ase de dados</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">a base escolhida foi a de </w:t>
      </w:r>
      <w:r>
        <w:t>diárias pagas aos servidores e colaboradores, pelos órgãos e entidades da Administração Pública Federal</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, no mês de dezembro do ano de 2017, disponível para download no site: </w:t>
      </w:r>
      <w:r>
        <w:t>http://www.portaldatransparencia.gov.br/downloads/mensal.asp?c=Diarias</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>U</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>suário</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>haverá o usuário DBA, que será o responsável pelo banco; também será disponibilizado o usuário aplicação sendo que este será o responsável pelo CRUD dos dados, e o usuário especializado o qual poderá acessar os relatórios</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">Ferramentas: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>SGBD</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>o</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> grupo optou pelo </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>MySQL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:t>por ser a linguage</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>m que lhe foi concedida para as atividades de mesa redondas, assim espera-se colocar os conhecimentos adquiridos em prol do trabalho final;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Linguagem de Programação:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">foi </w:t>
      </w:r>
      <w:r>
        <w:t>escolhido a linguagem orientada a objetos</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Java</w:t>
      </w:r>
      <w:r>
        <w:t>, pois o grupo já tem uma certa experiência com a linguagem</w:t>
      </w:r>
      <w:r>
        <w:t>, então</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> acredita</w:t>
      </w:r>
      <w:r>
        <w:t>mos</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">que tal escolha irá nos </w:t>
      </w:r>
      <w:r>
        <w:t>ajudar para a criação da aplicação.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">Ciência do Grupo: </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
